$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder text: "23/08/2022" -> "24/08/2022"
#    This placeholder lives on the Slide Master and on every slide layout
#    (the date auto-field that PowerPoint re-stamps on save).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "23/08/2022") {
                $shp.TextFrame.TextRange.Text = "24/08/2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

$customLayouts = $master.CustomLayouts
for ($li = 1; $li -le $customLayouts.Count; $li++) {
    $layout = $customLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Reposition "Picture 3" on slide 1.
#    Target EMU offset: x=8831737, y=5986015 (was x=9441337, y=5735459).
#    Shape.Left/Top are COM Singles (float32) truncated to EMU, so use the
#    float64 literal whose float32 cast lands exactly on the target EMU.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$pic = $slide1.Shapes.Item("Picture 3")
$pic.Left = 695.4124145507812
$pic.Top = 471.33978271484375

Write-Output "done"
